# Reorders data rows 2-19 (columns A:T) on the active sheet according to
# the mapping: destination row -> source row (content that should end up
# there). All rows are first snapshotted so the in-place permutation does
# not clobber values before they are copied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 19
$lastCol = 20   # column T

# destination row -> source row
$mapping = @{
    2  = 19
    3  = 4
    4  = 2
    5  = 13
    6  = 15
    7  = 17
    8  = 18
    9  = 5
    10 = 16
    11 = 12
    12 = 11
    13 = 8
    14 = 14
    15 = 10
    16 = 3
    17 = 7
    18 = 6
    19 = 9
}

# Snapshot current values (and number formats) for every row/column so the
# permutation can be applied without reading already-overwritten data.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $rowVals[$c] = $cell.Value2
    }
    $snapshot[$r] = $rowVals
}

# Write snapshotted values into their new (destination) rows.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $rowVals = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $rowVals[$c]
    }
}
